$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 43 (hunk 0)
$ws.Range("H43").Value = 6181.125
$ws.Range("I43").Value = 5271.143
$ws.Range("K43").Value = 5271.143
$ws.Range("M43").Value = -5202.143
# row 64 (hunk 1)
$ws.Range("H64").Value = 3112.8572
$ws.Range("J64").Value = 3112.8572
$ws.Range("L64").Value = 3112.8572
$ws.Range("N64").Value = -3608.8572
# row 67 (hunk 2)
$ws.Range("H67").Value = 3112.8572
$ws.Range("J67").Value = 3112.8572
$ws.Range("L67").Value = 3112.8572
$ws.Range("N67").Value = -4828.8572
# row 100 (hunk 3)
$ws.Range("H100").Value = 351
$ws.Range("I100").Value = 351
$ws.Range("K100").Value = 351
$ws.Range("M100").Value = 190
# row 113 (hunk 4)
$ws.Range("H113").Value = 2178
$ws.Range("I113").Value = 2360
$ws.Range("J113").Value = 1450
$ws.Range("K113").Value = 2360
$ws.Range("L113").Value = 1450
$ws.Range("M113").Value = 894
$ws.Range("N113").Value = -7958
# row 138 (hunk 5)
$ws.Range("H138").Value = 2454.3208
$ws.Range("I138").Value = 1853
$ws.Range("J138").Value = 3180.9167
$ws.Range("K138").Value = 5559
$ws.Range("L138").Value = 9542.750100000001
$ws.Range("M138").Value = -419
$ws.Range("N138").Value = -19822.7501

$ws = $wb.Worksheets.Item("ARM")
# row 61 (hunk 6)
$ws.Range("H61").Value = 2329.8333
$ws.Range("I61").Value = 2329.8333
$ws.Range("K61").Value = 2329.8333
$ws.Range("M61").Value = -2117.8333
# row 101 (hunk 7)
$ws.Range("H101").Value = 26801
$ws.Range("J101").Value = 26801
$ws.Range("L101").Value = 26801
$ws.Range("N101").Value = -33291
# row 105 (hunk 8)
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# row 122 (hunk 9)
$ws.Range("H122").Value = 928962.75
$ws.Range("I122").Value = 1269936.4
$ws.Range("J122").Value = 19699.666
$ws.Range("K122").Value = 3809809.2
$ws.Range("L122").Value = 59098.99800000001
$ws.Range("M122").Value = -3807359.2
$ws.Range("N122").Value = -63998.99800000001
# row 136 (hunk 10)
$ws.Range("H136").Value = 2329.8333
$ws.Range("I136").Value = 2329.8333
$ws.Range("K136").Value = 6989.499899999999
$ws.Range("M136").Value = -4439.499899999999

$ws = $wb.Worksheets.Item("BSM")
# row 107 (hunk 11)
$ws.Range("H107").Value = 1466.7368
$ws.Range("I107").Value = 1435.625
$ws.Range("K107").Value = 1435.625
$ws.Range("M107").Value = 484.375
# row 110 (hunk 12)
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180

$ws = $wb.Worksheets.Item("CRP")
# row 15 (hunk 13)
$ws.Range("H15").Value = 15282.333
$ws.Range("I15").Value = 14419.5
$ws.Range("J15").Value = 17008
$ws.Range("K15").Value = 14419.5
$ws.Range("L15").Value = 17008
$ws.Range("M15").Value = -14249.5
$ws.Range("N15").Value = -17348
# row 29 (hunk 14)
$ws.Range("H29").Value = 5412
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 5618
$ws.Range("K29").Value = 5000
$ws.Range("L29").Value = 5618
$ws.Range("M29").Value = -4707
$ws.Range("N29").Value = -6204
# row 31 (hunk 15)
$ws.Range("H31").Value = 3352.6667
$ws.Range("J31").Value = 7330
$ws.Range("L31").Value = 7330
$ws.Range("N31").Value = -7920
# row 34 (hunk 16)
$ws.Range("H34").Value = 3352.6667
$ws.Range("J34").Value = 7330
$ws.Range("L34").Value = 7330
$ws.Range("N34").Value = -7734
# row 60 (hunk 17)
$ws.Range("H60").Value = 13530.462
$ws.Range("J60").Value = 26658.666
$ws.Range("L60").Value = 26658.666
$ws.Range("N60").Value = -27680.666
# row 107 (hunk 18)
$ws.Range("H107").Value = 62500490
$ws.Range("I107").Value = 83333610
$ws.Range("J107").Value = 1124.5
$ws.Range("K107").Value = 83333610
$ws.Range("L107").Value = 1124.5
$ws.Range("M107").Value = -83331690
$ws.Range("N107").Value = -4964.5

$ws = $wb.Worksheets.Item("CUL")
# row 4 (hunk 19)
$ws.Range("H4").Value = 1764276.4
$ws.Range("J4").Value = 898.2
$ws.Range("L4").Value = 2694.6
$ws.Range("N4").Value = -2918.6
# row 37 (hunk 20)
$ws.Range("H37").Value = 99166.664
$ws.Range("J37").Value = 99166.664
$ws.Range("L37").Value = 297499.992
$ws.Range("N37").Value = -297723.992
# row 99 (hunk 21)
$ws.Range("H99").Value = 999
$ws.Range("I99").Value = 999
$ws.Range("K99").Value = 2997
$ws.Range("M99").Value = -751
# row 121 (hunk 22)
$ws.Range("H121").Value = 714.6667
$ws.Range("J121").Value = 800
$ws.Range("L121").Value = 2400
$ws.Range("N121").Value = -5020
# row 131 (hunk 23)
$ws.Range("H131").Value = 1250
$ws.Range("J131").Value = 1250
$ws.Range("L131").Value = 3750
$ws.Range("N131").Value = -13830

$ws = $wb.Worksheets.Item("GSM")
# row 18 (hunk 24)
$ws.Range("H18").Value = 3344533.2
$ws.Range("J18").Value = 16800
$ws.Range("L18").Value = 16800
$ws.Range("N18").Value = -17386
# row 100 (hunk 25)
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164
# row 101 (hunk 26)
$ws.Range("H101").Value = 46551.668
$ws.Range("J101").Value = 46551.668
$ws.Range("L101").Value = 46551.668
$ws.Range("N101").Value = -53041.668
# row 102 (hunk 27)
$ws.Range("H102").Value = 8326.532999999999
$ws.Range("I102").Value = 8564.143
$ws.Range("K102").Value = 8564.143
$ws.Range("M102").Value = -6942.143
# row 104 (hunk 28)
$ws.Range("H104").Value = 22613.5
$ws.Range("J104").Value = 22613.5
$ws.Range("L104").Value = 22613.5
$ws.Range("N104").Value = -29601.5
# row 122 (hunk 29)
$ws.Range("H122").Value = 144550.72
$ws.Range("I122").Value = 1976
$ws.Range("K122").Value = 5928
$ws.Range("M122").Value = -3478

$ws = $wb.Worksheets.Item("LTW")
# row 22 (hunk 30)
$ws.Range("H22").Value = 6034.75
$ws.Range("I22").Value = 7119.5
$ws.Range("J22").Value = 4950
$ws.Range("K22").Value = 7119.5
$ws.Range("L22").Value = 4950
$ws.Range("M22").Value = -6824.5
$ws.Range("N22").Value = -5540
# row 27 (hunk 31)
$ws.Range("H27").Value = 6034.75
$ws.Range("I27").Value = 7119.5
$ws.Range("J27").Value = 4950
$ws.Range("K27").Value = 7119.5
$ws.Range("L27").Value = 4950
$ws.Range("M27").Value = -7012.5
$ws.Range("N27").Value = -5164
# row 61 (hunk 32)
$ws.Range("H61").Value = 37038370
$ws.Range("I61").Value = 37038370
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 37038370
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -37038168
$ws.Range("N61").ClearContents()
# row 105 (hunk 33)
$ws.Range("H105").Value = 34807.5
$ws.Range("J105").Value = 34807.5
$ws.Range("L105").Value = 34807.5
$ws.Range("N105").Value = -41795.5
# row 113 (hunk 34)
$ws.Range("H113").Value = 37038370
$ws.Range("I113").Value = 37038370
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 37038370
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -37036200
$ws.Range("N113").ClearContents()
# row 132 (hunk 35)
$ws.Range("H132").Value = 2749.75
$ws.Range("I132").Value = 999.5
$ws.Range("K132").Value = 2998.5
$ws.Range("M132").Value = -468.5

$ws = $wb.Worksheets.Item("WVR")
# row 75 (hunk 36)
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# row 78 (hunk 37)
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# row 122 (hunk 38)
$ws.Range("H122").Value = 2551
$ws.Range("I122").Value = 2390
$ws.Range("K122").Value = 7170
$ws.Range("M122").Value = -4720
# row 132 (hunk 39)
$ws.Range("H132").Value = 2775
$ws.Range("I132").Value = 2069.3076
$ws.Range("K132").Value = 6207.9228
$ws.Range("M132").Value = -3677.9228
# row 136 (hunk 40)
$ws.Range("H136").Value = 1665.4
$ws.Range("I136").Value = 1081.75
$ws.Range("K136").Value = 3245.25
$ws.Range("M136").Value = -695.25
